$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.333.21"
$ws.Range("E2").Value = "  +5.56%  "
$ws.Range("D3").Value = "'3.005.25"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'581.74"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "'163.07"
$ws.Range("E6").Value = "  +11.80%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.519"
$ws.Range("E8").Value = "  +3.26%  "
$ws.Range("D9").Value = "'3.000.22"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("D10").Value = "'6.62"
$ws.Range("E10").Value = "  -4.90%  "
$ws.Range("D11").Value = "'0.155"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("E12").Value = "  +4.66%  "
$ws.Range("E13").Value = "  +5.29%  "
$ws.Range("D14").Value = "'34.82"
$ws.Range("E14").Value = "  +5.44%  "
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "'66.298.37"
$ws.Range("E16").Value = "  +5.83%  "
$ws.Range("D17").Value = "'3.502.26"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").Value = "'6.94"
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("D19").Value = "'3.010.18"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").Value = "'456.65"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("E21").Value = "  +5.14%  "
$ws.Range("D22").Value = "'0.689"
$ws.Range("E22").Value = "  +3.66%  "
$ws.Range("D23").Value = "'7.37"
$ws.Range("E23").Value = "  +6.40%  "
$ws.Range("D24").Value = "'82.36"
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("E25").Value = "  +13.04%  "
$ws.Range("D26").Value = "'12.38"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").Value = "'10.50"
$ws.Range("E27").Value = "  +3.92%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +16.24%  "
$ws.Range("E30").Value = "  +18.40%  "
$ws.Range("D31").Value = "'0.0000105"
$ws.Range("E31").Value = "  -7.01%  "
$ws.Range("D32").Value = "'2.62"
$ws.Range("E32").Value = "  +4.36%  "
$ws.Range("D33").Value = "'27.28"
$ws.Range("E33").Value = "  +4.93%  "
$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'5.89"
$ws.Range("E36").Value = "  +8.30%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "'0.991"
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("D38").Value = "'2.18"
$ws.Range("E38").Value = "  +13.09%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'49.89"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("E41").Value = "  +15.24%  "
$ws.Range("D42").Value = "'0.123"
$ws.Range("E42").Value = "  +7.41%  "
$ws.Range("D43").Value = "'44.16"
$ws.Range("E43").Value = "  +7.03%  "
$ws.Range("D44").Value = "'8.44"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("D45").Value = "'401.88"
$ws.Range("E45").Value = "  +13.41%  "
$ws.Range("D46").Value = "'0.0362"
$ws.Range("E46").Value = "  +6.11%  "
$ws.Range("D47").Value = "'2.798.05"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D50").Value = "'24.09"
$ws.Range("E50").Value = "  +11.06%  "
$ws.Range("E51").Value = "  +4.16%  "
